$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their text representation (e.g. "1.000",
# "42.00") instead of Excel auto-converting the numeric-looking strings to numbers
# and stripping "insignificant" trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.847.26'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '1.755.74'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('D5').Value = '327.46'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '0.4652'
$ws.Range('E7').Value = '  +1.61%  '
$ws.Range('D8').Value = '0.3494'
$ws.Range('E8').Value = '  -1.96%  '
$ws.Range('D9').Value = '42.00'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').Value = '0.07356'
$ws.Range('E10').Value = '  -1.43%  '
$ws.Range('D11').Value = '1.080'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = '20.49'
$ws.Range('E13').Value = '  -1.54%  '
$ws.Range('D14').Value = '5.979'
$ws.Range('D15').Value = '7.142'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').Value = '1.756.31'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '92.06'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').Value = '0.00001054'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').Value = '0.06401'
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('E21').Value = '  -1.66%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '27.860.28'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D24').Value = '11.13'
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('D25').Value = '2.150'
$ws.Range('E25').Value = '  +3.83%  '
$ws.Range('D26').Value = '161.50'
$ws.Range('E26').Value = '  -2.41%  '
$ws.Range('D27').Value = '20.00'
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').Value = '1.956.99'
$ws.Range('D29').Value = '2.150'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').Value = '122.68'
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('D31').Value = '1.069'
$ws.Range('E31').Value = '  -0.89%  '
$ws.Range('D32').Value = '0.09299'
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('D33').Value = '3.651'
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('D34').Value = '5.540'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.02263'
$ws.Range('E35').Value = '  -0.66%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = '11.62'
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('D37').Value = '0.06064'
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('D38').Value = '0.2061'
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('D39').Value = '4.894'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').Value = '0.6134'
$ws.Range('E40').Value = '  -2.16%  '
$ws.Range('D41').Value = '1.175'
$ws.Range('D42').Value = '7.759'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').Value = '1.351'
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('E44').Value = '  -1.96%  '
$ws.Range('D45').Value = '3.732'
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('D46').Value = '0.5770'
$ws.Range('E46').Value = '  -1.49%  '
$ws.Range('D47').Value = '122.74'
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('D48').Value = '1.921'
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('D49').Value = '0.06795'
$ws.Range('E49').Value = '  -1.69%  '
$ws.Range('D50').Value = '1.119'
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('D51').Value = '72.09'
$ws.Range('E51').Value = '  -0.04%  '
